$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VO IDs range")

# Update the reserved ID range for "vaccine adjuvant" to reflect the two newly
# assigned IDs (VO_0005510 and VO_0005511 are now used), issues #759 and #760.
$ws.Range("A2").Value = "VO_0005512 - VO_0005560"

# Move the active selection, matching where the author left off editing.
$ws.Range("D5").Select()
